# Refresh the crypto price/volume table with the latest scraped values.
# D-column prices that look like plain numbers are written with a leading
# apostrophe (forcing text) and the cell style is reset to "Normal" right
# after, so the stored value stays a literal string (matching the
# dotted/locale price formatting already used in the sheet, e.g.
# "69.277.98") instead of being auto-converted into a float by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.277.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").Value = "'3.409.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'582.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").Value = "'178.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +7.94%  "

$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("D11").Value = "'48.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D12").Value = "'0.0000282"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.24%  "

$ws.Range("D13").Value = "'681.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "

$ws.Range("D14").Value = "'8.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.28%  "

$ws.Range("D15").Value = "'3.953.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("D16").Value = "'69.392.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("D17").Value = "'3.405.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("D20").Value = "'11.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").Value = "'0.912"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.87%  "

$ws.Range("D22").Value = "'5.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.59%  "

$ws.Range("D23").Value = "'17.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.67%  "

$ws.Range("D24").Value = "'100.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.65%  "

$ws.Range("E25").Value = "  -0.34%  "

$ws.Range("D27").Value = "'9.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.01%  "

$ws.Range("D28").Value = "'33.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("D29").Value = "'8.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.64%  "

$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").Value = "'3.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.88%  "

$ws.Range("D32").Value = "'557.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.95%  "

$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").Value = "'57.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").Value = "'3.613.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.77%  "

$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").Value = "'35.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("D40").Value = "'0.0₃0744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.53%  "

$ws.Range("D41").Value = "'3.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.22%  "

$ws.Range("E42").Value = "  +3.30%  "

$ws.Range("D43").Value = "'3.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.61%  "

$ws.Range("D44").Value = "'0.0425"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.28%  "

$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("E46").Value = "  +1.02%  "

$ws.Range("D48").Value = "'1.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.02%  "

$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").Value = "'131.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("E51").Value = "  +3.37%  "
